$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7771000862121582
$ws.Range("B1").Value = 0.8826313018798828
$ws.Range("C1").Value = 3.512059688568115
$ws.Range("D1").Value = 2.117461681365967
$ws.Range("E1").Value = 1.041028141975403
